$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column O (Clinic Notes), shifting existing
# columns O:AC to P:AD.
$ws.Columns("O:O").Insert()

# New header cell for the inserted column.
$ws.Range("O1").Value = "Birth Year"

# Match the new column's width to its left neighbor (Age Units), mirroring
# Excel's default "insert column" formatting behavior.
$ws.Columns("O:O").ColumnWidth = $ws.Columns("N:N").ColumnWidth

# Update the active selection to match the saved view state.
$ws.Range("N3").Select()
